$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2898.1333
$ws.Range("I4").Value = 3398.1667
$ws.Range("J4").Value = 2564.7778
$ws.Range("K4").Value = 3398.1667
$ws.Range("L4").Value = 2564.7778
$ws.Range("M4").Value = -3284.1667
$ws.Range("N4").Value = -2792.7778

$ws.Range("H55").Value = 1007.4
$ws.Range("I55").Value = 550
$ws.Range("J55").Value = 1312.3334
$ws.Range("K55").Value = 550
$ws.Range("L55").Value = 1312.3334
$ws.Range("M55").Value = -336
$ws.Range("N55").Value = -1740.3334

$ws.Range("H94").Value = 639.4286
$ws.Range("I94").Value = 899.25
$ws.Range("J94").Value = 293
$ws.Range("K94").Value = 899.25
$ws.Range("L94").Value = 293
$ws.Range("M94").Value = -448.25
$ws.Range("N94").Value = -1195

$ws.Range("H100").Value = 3701
$ws.Range("I100").Value = 2836.3333
$ws.Range("K100").Value = 2836.3333
$ws.Range("M100").Value = -2295.3333

$ws.Range("H107").Value = 807.9167
$ws.Range("I107").Value = 822.65216
$ws.Range("J107").Value = 469
$ws.Range("K107").Value = 822.65216
$ws.Range("L107").Value = 469
$ws.Range("M107").Value = 1097.34784
$ws.Range("N107").Value = -4309

$ws.Range("H112").Value = 9132.333000000001
$ws.Range("J112").Value = 10560.7
$ws.Range("L112").Value = 31682.1
$ws.Range("N112").Value = -33898.10000000001

$ws.Range("H137").Value = 2394.4
$ws.Range("I137").Value = 1951.2
$ws.Range("K137").Value = 5853.6
$ws.Range("M137").Value = -3303.6

$ws.Range("H141").Value = 2187.125
$ws.Range("I141").Value = 1818.8572
$ws.Range("J141").Value = 4765
$ws.Range("K141").Value = 5456.571599999999
$ws.Range("L141").Value = 14295
$ws.Range("M141").Value = -276.5715999999993
$ws.Range("N141").Value = -24655

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9646.637000000001
$ws.Range("I32").Value = 6957.2676
$ws.Range("K32").Value = 6957.2676
$ws.Range("M32").Value = -6670.2676

$ws.Range("H74").Value = 4098.2666
$ws.Range("I74").Value = 2035.6154
$ws.Range("J74").Value = 17505.5
$ws.Range("K74").Value = 2035.6154
$ws.Range("L74").Value = 17505.5
$ws.Range("M74").Value = -1161.6154
$ws.Range("N74").Value = -19253.5

$ws.Range("H77").Value = 4098.2666
$ws.Range("I77").Value = 2035.6154
$ws.Range("J77").Value = 17505.5
$ws.Range("K77").Value = 10178.077
$ws.Range("L77").Value = 87527.5
$ws.Range("M77").Value = -5810.076999999999
$ws.Range("N77").Value = -96263.5

$ws.Range("H109").Value = 65999
$ws.Range("J109").Value = 65999
$ws.Range("L109").Value = 65999
$ws.Range("N109").Value = -68773

$ws.Range("H122").Value = 6166.647
$ws.Range("I122").Value = 4822.3335
$ws.Range("K122").Value = 14467.0005
$ws.Range("M122").Value = -12017.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 30907
$ws.Range("J58").Value = 29361.25
$ws.Range("L58").Value = 29361.25
$ws.Range("N58").Value = -29949.25

$ws.Range("H126").Value = 70000
$ws.Range("J126").Value = 70000
$ws.Range("L126").Value = 70000
$ws.Range("N126").Value = -79880

$ws.Range("H134").Value = 4168
$ws.Range("I134").Value = 3290.7812
$ws.Range("J134").Value = 5727.5
$ws.Range("K134").Value = 9872.3436
$ws.Range("L134").Value = 17182.5
$ws.Range("M134").Value = -7337.3436
$ws.Range("N134").Value = -22252.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7626.718
$ws.Range("I31").Value = 3240.8845
$ws.Range("J31").Value = 16398.385
$ws.Range("K31").Value = 3240.8845
$ws.Range("L31").Value = 16398.385
$ws.Range("M31").Value = -2945.8845
$ws.Range("N31").Value = -16988.385

$ws.Range("H34").Value = 7626.718
$ws.Range("I34").Value = 3240.8845
$ws.Range("J34").Value = 16398.385
$ws.Range("K34").Value = 3240.8845
$ws.Range("L34").Value = 16398.385
$ws.Range("M34").Value = -3038.8845
$ws.Range("N34").Value = -16802.385

$ws.Range("H44").Value = 60000
$ws.Range("J44").Value = 60000
$ws.Range("L44").Value = 60000
$ws.Range("N44").Value = -60884

$ws.Range("H52").Value = 134000
$ws.Range("I52").Value = 134000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 134000
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -133706
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5630.95
$ws.Range("J131").Value = 6795.2144
$ws.Range("L131").Value = 20385.6432
$ws.Range("N131").Value = -30465.6432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4630.6206
$ws.Range("J122").Value = 15525
$ws.Range("L122").Value = 46575
$ws.Range("N122").Value = -51475

$ws.Range("H132").Value = 5197.58
$ws.Range("I132").Value = 4402.4
$ws.Range("K132").Value = 13207.2
$ws.Range("M132").Value = -10677.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2177
$ws.Range("J22").Value = 4722
$ws.Range("L22").Value = 4722
$ws.Range("N22").Value = -5312

$ws.Range("H27").Value = 2177
$ws.Range("J27").Value = 4722
$ws.Range("L27").Value = 4722
$ws.Range("N27").Value = -4936

$ws.Range("H122").Value = 3748.3333
$ws.Range("I122").Value = 2236.762
$ws.Range("J122").Value = 7275.3335
$ws.Range("K122").Value = 6710.286
$ws.Range("L122").Value = 21826.0005
$ws.Range("M122").Value = -4260.286
$ws.Range("N122").Value = -26726.0005

$ws.Range("H132").Value = 7175.3145
$ws.Range("I132").Value = 6469.64
$ws.Range("J132").Value = 8939.5
$ws.Range("K132").Value = 19408.92
$ws.Range("L132").Value = 26818.5
$ws.Range("M132").Value = -16878.92
$ws.Range("N132").Value = -31878.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 259500.75
$ws.Range("I3").Value = 6003
$ws.Range("J3").Value = 344000
$ws.Range("K3").Value = 6003
$ws.Range("L3").Value = 344000
$ws.Range("M3").Value = -5889
$ws.Range("N3").Value = -344228

$ws.Range("H62").Value = 7664.8335
$ws.Range("I62").Value = 3999.5
$ws.Range("J62").Value = 9497.5
$ws.Range("K62").Value = 3999.5
$ws.Range("L62").Value = 9497.5
$ws.Range("M62").Value = -3375.5
$ws.Range("N62").Value = -10745.5

$ws.Range("H65").Value = 7664.8335
$ws.Range("I65").Value = 3999.5
$ws.Range("J65").Value = 9497.5
$ws.Range("K65").Value = 19997.5
$ws.Range("L65").Value = 47487.5
$ws.Range("M65").Value = -16877.5
$ws.Range("N65").Value = -53727.5

$ws.Range("H107").Value = 3360.9524
$ws.Range("I107").Value = 2938.3333
$ws.Range("J107").Value = 4417.5
$ws.Range("K107").Value = 8814.999899999999
$ws.Range("L107").Value = 13252.5
$ws.Range("M107").Value = -6894.999899999999
$ws.Range("N107").Value = -17092.5

$ws.Range("H122").Value = 5272.263
$ws.Range("I122").Value = 5274.8823
$ws.Range("K122").Value = 15824.6469
$ws.Range("M122").Value = -13374.6469

$ws.Range("H136").Value = 6149.1577
$ws.Range("I136").Value = 5071.5835
$ws.Range("J136").Value = 7996.4287
$ws.Range("K136").Value = 15214.7505
$ws.Range("L136").Value = 23989.2861
$ws.Range("M136").Value = -12664.7505
$ws.Range("N136").Value = -29089.2861
